# Add season-record columns (Wins / Losses / Ties) to the roster table.
#
# The existing sheet holds one team's roster (header row 1, data rows
# 2-51, columns A:AC). We append three new columns — AD "Wins",
# AE "Losses", AF "Ties" — carrying the team's 1998 season record
# (54-108-0) repeated for every player row, matching the header's bold /
# centered / bordered style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 51

# Clone the header formatting (bold font, centered/top alignment, thin
# border) from an existing header cell onto the three new header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record for every player on the roster.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 54   # AD: Wins
    $ws.Cells.Item($r, 31).Value = 108  # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF: Ties
}
